$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 536; this shifts the existing rows
# 536-601 down to 537-602 (matches the xlsx diff: dimension grows from
# A1:R601 to A1:R602, and every row from the old 536 onward is now one
# row lower).
$ws.Rows.Item(536).Insert()

# Populate the newly inserted row 536 with the new record's data. The
# columns that are constant across the whole sheet (A, B, C, E, F, G, H,
# I, N, Q, R) keep the same values as their neighbouring rows.
$ws.Range("A536").Value = 3
$ws.Range("B536").Value = "Femacal de La Calera"
$ws.Range("C536").Value = "Coquimbo"
$ws.Range("D536").Value = 45142
$ws.Range("E536").Value = 5
$ws.Range("F536").Value = 100114013
$ws.Range("G536").Value = "Zanahoria"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Primera"
$ws.Range("J536").Value = 120
$ws.Range("K536").Value = 7500
$ws.Range("L536").Value = 7500
$ws.Range("M536").Value = 7500
$ws.Range("N536").Value = "$/saco 20 kilos"
$ws.Range("O536").Value = "Provincia de Quillota"
$ws.Range("P536").Value = 375
$ws.Range("Q536").Value = 20
$ws.Range("R536").Value = "Hortaliza"
